$d = $word.ActiveDocument

function ReplaceExact($find) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $find, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# --- Merge the runs that were artificially split around "(" ")" with gramStart/gramEnd proofErr marks ---
ReplaceExact("__construct()")
ReplaceExact("getidTarea()")
ReplaceExact("getTarea()")
ReplaceExact("getidUsuario()")
ReplaceExact("getUsuario_Tweeter()")
ReplaceExact("getHastag()")
ReplaceExact("getDia_inicio()")
ReplaceExact("getDia_Fin()")
ReplaceExact("getHora_Inicio()")
ReplaceExact("getHora_Fin()")

# --- Merge "twe" + bookmark("_GoBack") + "eter" into a single "tweeter" run ---
ReplaceExact("tweeter")

# --- Fix the typo: "hasta" -> "hashtag" (the actual content fix from the commit message) ---
$rng1 = $d.Content
$rng1.Find.Execute("asignar el hasta de la tarea", $true, $false, $false, $false, $false, $true, 1, $false, "asignar el hashtag de la tarea", 2)
$rng2 = $d.Content
$rng2.Find.Execute("obtener hasta de la tarea", $true, $false, $false, $false, $false, $true, 1, $false, "obtener hashtag de la tarea", 2)
